# "Cambios en la rama1 pero no en master"
#
# Target change (per the OOXML diff):
#   - paragraph "dfsfd" becomes "Dfsfd" (split "D" / "fsfd", wrapped in
#     spellcheck <w:proofErr> start/end marks)
#   - a brand-new paragraph "Cambio para ver con git diff" is inserted
#     right after it (with "git" and "diff" individually wrapped in
#     <w:proofErr> marks, as Word's spell checker would do while typing)
#   - the trailing _GoBack bookmark, which used to sit at the end of the
#     "dfsfd" paragraph, ends up alone in its own (new, empty) paragraph
#     after the inserted text.
#
# Word's Range.InsertXML() happens to replace far more than the target
# range when the injected fragment contains more than one <w:p>, so we
# first grow the document to the right number of paragraphs with plain
# InsertParagraphAfter() calls, and then stamp the exact OOXML for each
# paragraph one at a time (one paragraph of payload per one paragraph of
# range -- that form is applied in place, as desired).

$d = $word.ActiveDocument

# Find the paragraph that holds "dfsfd" (index 2 in the original file).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "dfsfd") {
        $target = $i
        break
    }
}
if ($target -eq $null) {
    throw "Could not locate the 'dfsfd' paragraph"
}

$wmain = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Make room: split off two new (still empty) paragraphs right after the
# "dfsfd" paragraph -- one will become the "Cambio..." paragraph, the
# other will become the bookmark-only paragraph that used to trail
# "dfsfd".
$anchor = $d.Paragraphs($target)
$anchor.Range.InsertParagraphAfter()
$anchor.Range.InsertParagraphAfter()

# 1) "dfsfd" -> "Dfsfd" (two runs, wrapped in proofErr spellStart/spellEnd)
$xmlDfsfd = '<w:document ' + $wmain + '><w:body>' +
    '<w:p>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>D</w:t></w:r>' +
        '<w:r><w:t>fsfd</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>' +
    '</w:body></w:document>'
$d.Paragraphs($target).Range.InsertXML($xmlDfsfd)

# 2) New paragraph: "Cambio para ver con git diff"
$xmlCambio = '<w:document ' + $wmain + '><w:body>' +
    '<w:p>' +
        '<w:r><w:t xml:space="preserve">Cambio para ver con </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>git</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>diff</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>' +
    '</w:body></w:document>'
$d.Paragraphs($target + 1).Range.InsertXML($xmlCambio)

# 3) Trailing paragraph keeps only the _GoBack bookmark.
$xmlBookmark = '<w:document ' + $wmain + '><w:body>' +
    '<w:p>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>' +
    '</w:body></w:document>'
$d.Paragraphs($target + 2).Range.InsertXML($xmlBookmark)

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output ($i.ToString() + ": [" + $d.Paragraphs($i).Range.Text + "]")
}
